# Auto-generated: apply the Ramuh_Profits (per-sheet leve data) value updates
# from the scheduled-runner commit. Only literal H..N value cells change;
# no formulas are present in these tables, so each target cell is written directly.

$wb = $excel.ActiveWorkbook

# ----- ALC sheet -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 93.818184
$ws.Range("I33").Value = 66
$ws.Range("J33").Value = 104.25
$ws.Range("K33").Value = 66
$ws.Range("L33").Value = 104.25
$ws.Range("M33").Value = 163
$ws.Range("N33").Value = -562.25
$ws.Range("H100").Value = 2390.0625
$ws.Range("I100").Value = 1689.375
$ws.Range("J100").Value = 3090.75
$ws.Range("K100").Value = 1689.375
$ws.Range("L100").Value = 3090.75
$ws.Range("M100").Value = -1148.375
$ws.Range("N100").Value = -4172.75
$ws.Range("H113").Value = 2623.2222
$ws.Range("I113").Value = 1488.75
$ws.Range("J113").Value = 3530.8
$ws.Range("K113").Value = 1488.75
$ws.Range("L113").Value = 3530.8
$ws.Range("M113").Value = 1765.25
$ws.Range("N113").Value = -10038.8
$ws.Range("H116").Value = 2665.6667
$ws.Range("I116").Value = 3701.25
$ws.Range("J116").Value = 2289.0908
$ws.Range("K116").Value = 3701.25
$ws.Range("L116").Value = 2289.0908
$ws.Range("M116").Value = -259.25
$ws.Range("N116").Value = -9173.0908
$ws.Range("H129").Value = 710
$ws.Range("I129").Value = 565
$ws.Range("K129").Value = 1695
$ws.Range("M129").Value = 3305
$ws.Range("H141").Value = 3567.8723
$ws.Range("I141").Value = 1242.7142
$ws.Range("J141").Value = 10349.583
$ws.Range("K141").Value = 3728.1426
$ws.Range("L141").Value = 31048.749
$ws.Range("M141").Value = 1451.8574
$ws.Range("N141").Value = -41408.749

# ----- ARM sheet -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1164.9032
$ws.Range("I45").Value = 1027.2307
$ws.Range("J45").Value = 1880.8
$ws.Range("K45").Value = 1027.2307
$ws.Range("L45").Value = 1880.8
$ws.Range("M45").Value = -650.2307000000001
$ws.Range("N45").Value = -2634.8
$ws.Range("H110").Value = 2247.1428
$ws.Range("I110").Value = 2600
$ws.Range("J110").Value = 1776.6666
$ws.Range("K110").Value = 2600
$ws.Range("L110").Value = 1776.6666
$ws.Range("M110").Value = -555
$ws.Range("N110").Value = -5866.6666
$ws.Range("H135").Value = 18833.334
$ws.Range("J135").Value = 18833.334
$ws.Range("L135").Value = 18833.334
$ws.Range("N135").Value = -28973.334

# ----- BSM sheet -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3483.8462
$ws.Range("I86").Value = 3784.3333
$ws.Range("J86").Value = 3226.2856
$ws.Range("K86").Value = 3784.3333
$ws.Range("L86").Value = 3226.2856
$ws.Range("M86").Value = -2661.3333
$ws.Range("N86").Value = -5472.2856
$ws.Range("H89").Value = 3483.8462
$ws.Range("I89").Value = 3784.3333
$ws.Range("J89").Value = 3226.2856
$ws.Range("K89").Value = 18921.6665
$ws.Range("L89").Value = 16131.428
$ws.Range("M89").Value = -13305.6665
$ws.Range("N89").Value = -27363.428
$ws.Range("H107").Value = 1794.7931
$ws.Range("I107").Value = 861.1177
$ws.Range("J107").Value = 3117.5
$ws.Range("K107").Value = 861.1177
$ws.Range("L107").Value = 3117.5
$ws.Range("M107").Value = 1058.8823
$ws.Range("N107").Value = -6957.5

# ----- CRP sheet -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3060
$ws.Range("I16").Value = 1360
$ws.Range("J16").Value = 4760
$ws.Range("K16").Value = 1360
$ws.Range("L16").Value = 4760
$ws.Range("M16").Value = -1073
$ws.Range("N16").Value = -5334
$ws.Range("H31").Value = 32911.5
$ws.Range("I31").Value = 1034.6818
$ws.Range("J31").Value = 91352.336
$ws.Range("K31").Value = 1034.6818
$ws.Range("L31").Value = 91352.336
$ws.Range("M31").Value = -739.6818000000001
$ws.Range("N31").Value = -91942.336
$ws.Range("H34").Value = 32911.5
$ws.Range("I34").Value = 1034.6818
$ws.Range("J34").Value = 91352.336
$ws.Range("K34").Value = 1034.6818
$ws.Range("L34").Value = 91352.336
$ws.Range("M34").Value = -832.6818000000001
$ws.Range("N34").Value = -91756.336
$ws.Range("H107").Value = 63511.312
$ws.Range("I107").Value = 112155.555
$ws.Range("J107").Value = 968.7143
$ws.Range("K107").Value = 112155.555
$ws.Range("L107").Value = 968.7143
$ws.Range("M107").Value = -110235.555
$ws.Range("N107").Value = -4808.7143
$ws.Range("H108").Value = 35342
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 35342
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 35342
$ws.Range("M108").ClearContents()
$ws.Range("N108").Value = -43022
$ws.Range("H109").Value = 19000
$ws.Range("J109").Value = 19000
$ws.Range("L109").Value = 19000
$ws.Range("N109").Value = -21080
$ws.Range("H110").Value = 42851
$ws.Range("J110").Value = 42851
$ws.Range("L110").Value = 42851
$ws.Range("N110").Value = -51031
$ws.Range("H112").Value = 29800.666
$ws.Range("J112").Value = 29800.666
$ws.Range("L112").Value = 29800.666
$ws.Range("N112").Value = -32754.666
$ws.Range("H113").Value = 3060
$ws.Range("I113").Value = 1360
$ws.Range("J113").Value = 4760
$ws.Range("K113").Value = 1360
$ws.Range("L113").Value = 4760
$ws.Range("M113").Value = 810
$ws.Range("N113").Value = -9100
$ws.Range("H114").Value = 28400
$ws.Range("J114").Value = 28400
$ws.Range("L114").Value = 28400
$ws.Range("N114").Value = -37078
$ws.Range("H116").Value = 48000
$ws.Range("J116").Value = 48000
$ws.Range("L116").Value = 48000
$ws.Range("N116").Value = -57178
$ws.Range("H117").Value = 46387.5
$ws.Range("J117").Value = 46387.5
$ws.Range("L117").Value = 46387.5
$ws.Range("N117").Value = -55565.5
$ws.Range("H118").Value = 37950
$ws.Range("J118").Value = 37950
$ws.Range("L118").Value = 37950
$ws.Range("N118").Value = -41264
$ws.Range("H119").Value = 48000
$ws.Range("J119").Value = 48000
$ws.Range("L119").Value = 48000
$ws.Range("N119").Value = -57676
$ws.Range("H121").Value = 48000
$ws.Range("J121").Value = 48000
$ws.Range("L121").Value = 48000
$ws.Range("N121").Value = -50620

# ----- CUL sheet -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 7123
$ws.Range("J121").Value = 7747.6562
$ws.Range("L121").Value = 23242.9686
$ws.Range("N121").Value = -25862.9686

# ----- GSM sheet -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3431.9211
$ws.Range("I102").Value = 3506.1428
$ws.Range("J102").Value = 3224.1
$ws.Range("K102").Value = 3506.1428
$ws.Range("L102").Value = 3224.1
$ws.Range("M102").Value = -1884.1428
$ws.Range("N102").Value = -6468.1
$ws.Range("H107").Value = 178.81818
$ws.Range("I107").Value = 173.76471
$ws.Range("J107").Value = 196
$ws.Range("K107").Value = 173.76471
$ws.Range("L107").Value = 196
$ws.Range("M107").Value = 1746.23529
$ws.Range("N107").Value = -4036

# ----- WVR sheet -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 251.8
$ws.Range("I107").Value = 241.21428
$ws.Range("K107").Value = 723.64284
$ws.Range("M107").Value = 1196.35716
